# Auto-generated edit script applying the Shiva_Profits market-data refresh
# (see commit: "chore: update Sheets via scheduled runner")
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 488080.44   # H15: 4993.095 -> 488080.44
$ws.Cells.Item(15, 9).Value = 488080.44   # I15: 4993.095 -> 488080.44
$ws.Cells.Item(15, 11).Value = 1464241.32   # K15: 14979.285 -> 1464241.32
$ws.Cells.Item(15, 13).Value = -1464072.32   # M15: -14810.285 -> -1464072.32
$ws.Cells.Item(19, 8).Value = 8019.7896   # H19: 11572.077 -> 8019.7896
$ws.Cells.Item(19, 9).Value = 872.2   # I19: 990.3333 -> 872.2
$ws.Cells.Item(19, 10).Value = 10572.5   # J19: 14746.6 -> 10572.5
$ws.Cells.Item(19, 11).Value = 872.2   # K19: 990.3333 -> 872.2
$ws.Cells.Item(19, 12).Value = 10572.5   # L19: 14746.6 -> 10572.5
$ws.Cells.Item(19, 13).Value = -697.2   # M19: -815.3333 -> -697.2
$ws.Cells.Item(19, 14).Value = -10922.5   # N19: -15096.6 -> -10922.5
$ws.Cells.Item(98, 8).Value = 1655.7222   # H98: 1729.8823 -> 1655.7222
$ws.Cells.Item(98, 9).Value = 1655.7222   # I98: 1729.8823 -> 1655.7222
$ws.Cells.Item(98, 11).Value = 1655.7222   # K98: 1729.8823 -> 1655.7222
$ws.Cells.Item(98, 13).Value = -157.7221999999999   # M98: -231.8823 -> -157.7221999999999
$ws.Cells.Item(112, 8).Value = 2996.7368   # H112: 2098.2856 -> 2996.7368
$ws.Cells.Item(112, 9).Value = 1899.75   # I112: 1900 -> 1899.75
$ws.Cells.Item(112, 10).Value = 3289.2666   # J112: 3288 -> 3289.2666
$ws.Cells.Item(112, 11).Value = 5699.25   # K112: 5700 -> 5699.25
$ws.Cells.Item(112, 12).Value = 9867.799800000001   # L112: 9864 -> 9867.799800000001
$ws.Cells.Item(112, 13).Value = -4591.25   # M112: -4592 -> -4591.25
$ws.Cells.Item(112, 14).Value = -12083.7998   # N112: -12080 -> -12083.7998
$ws.Cells.Item(122, 8).Value = 1655.7222   # H122: 1729.8823 -> 1655.7222
$ws.Cells.Item(122, 9).Value = 1655.7222   # I122: 1729.8823 -> 1655.7222
$ws.Cells.Item(122, 11).Value = 4967.1666   # K122: 5189.6469 -> 4967.1666
$ws.Cells.Item(122, 13).Value = -2517.1666   # M122: -2739.6469 -> -2517.1666
$ws.Cells.Item(132, 8).Value = 6319.0586   # H132: 6445.94 -> 6319.0586
$ws.Cells.Item(132, 9).Value = 5000.1055   # I132: 5135.919 -> 5000.1055
$ws.Cells.Item(132, 11).Value = 15000.3165   # K132: 15407.757 -> 15000.3165
$ws.Cells.Item(132, 13).Value = -12470.3165   # M132: -12877.757 -> -12470.3165
$ws.Cells.Item(135, 8).Value = 12822043   # H135: 13159805 -> 12822043
$ws.Cells.Item(135, 9).Value = 17858354   # I135: 17858464 -> 17858354
$ws.Cells.Item(135, 10).Value = 2339.2727   # J135: 3561.8 -> 2339.2727
$ws.Cells.Item(135, 11).Value = 160725186   # K135: 160726176 -> 160725186
$ws.Cells.Item(135, 12).Value = 21053.4543   # L135: 32056.2 -> 21053.4543
$ws.Cells.Item(135, 13).Value = -160722651   # M135: -160723641 -> -160722651
$ws.Cells.Item(135, 14).Value = -26123.4543   # N135: -37126.2 -> -26123.4543
$ws.Cells.Item(138, 8).Value = 213341810   # H138: 266675440 -> 213341810
$ws.Cells.Item(138, 9).Value = 333340800   # I138: 500007460 -> 333340800
$ws.Cells.Item(138, 10).Value = 33343332   # J138: 33343440 -> 33343332
$ws.Cells.Item(138, 11).Value = 1000022400   # K138: 1500022380 -> 1000022400
$ws.Cells.Item(138, 12).Value = 100029996   # L138: 100030320 -> 100029996
$ws.Cells.Item(138, 13).Value = -1000017260   # M138: -1500017240 -> -1000017260
$ws.Cells.Item(138, 14).Value = -100040276   # N138: -100040600 -> -100040276
$ws.Cells.Item(140, 8).Value = 67917.5   # H140: 67972.75 -> 67917.5
$ws.Cells.Item(140, 10).Value = 67917.5   # J140: 67972.75 -> 67917.5
$ws.Cells.Item(140, 12).Value = 67917.5   # L140: 67972.75 -> 67917.5
$ws.Cells.Item(140, 14).Value = -78277.5   # N140: -78332.75 -> -78277.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4668.22   # H32: 3526.0376 -> 4668.22
$ws.Cells.Item(32, 9).Value = 3735.9272   # I32: 2978.203 -> 3735.9272
$ws.Cells.Item(32, 10).Value = 17487.25   # J32: 6962.4546 -> 17487.25
$ws.Cells.Item(32, 11).Value = 3735.9272   # K32: 2978.203 -> 3735.9272
$ws.Cells.Item(32, 12).Value = 17487.25   # L32: 6962.4546 -> 17487.25
$ws.Cells.Item(32, 13).Value = -3448.9272   # M32: -2691.203 -> -3448.9272
$ws.Cells.Item(32, 14).Value = -18061.25   # N32: -7536.4546 -> -18061.25
$ws.Cells.Item(61, 8).Value = 2474.074   # H61: 5134.533 -> 2474.074
$ws.Cells.Item(61, 9).Value = 1791   # I61: 2746.5557 -> 1791
$ws.Cells.Item(61, 10).Value = 4864.8335   # J61: 8716.5 -> 4864.8335
$ws.Cells.Item(61, 11).Value = 1791   # K61: 2746.5557 -> 1791
$ws.Cells.Item(61, 12).Value = 4864.8335   # L61: 8716.5 -> 4864.8335
$ws.Cells.Item(61, 13).Value = -1579   # M61: -2534.5557 -> -1579
$ws.Cells.Item(61, 14).Value = -5288.8335   # N61: -9140.5 -> -5288.8335
$ws.Cells.Item(132, 8).Value = 2354.5625   # H132: 2275.4707 -> 2354.5625
$ws.Cells.Item(132, 9).Value = 1727.3572   # I132: 1679.5333 -> 1727.3572
$ws.Cells.Item(132, 11).Value = 5182.071599999999   # K132: 5038.5999 -> 5182.071599999999
$ws.Cells.Item(132, 13).Value = -2652.071599999999   # M132: -2508.5999 -> -2652.071599999999
$ws.Cells.Item(136, 8).Value = 2474.074   # H136: 5134.533 -> 2474.074
$ws.Cells.Item(136, 9).Value = 1791   # I136: 2746.5557 -> 1791
$ws.Cells.Item(136, 10).Value = 4864.8335   # J136: 8716.5 -> 4864.8335
$ws.Cells.Item(136, 11).Value = 5373   # K136: 8239.667099999999 -> 5373
$ws.Cells.Item(136, 12).Value = 14594.5005   # L136: 26149.5 -> 14594.5005
$ws.Cells.Item(136, 13).Value = -2823   # M136: -5689.667099999999 -> -2823
$ws.Cells.Item(136, 14).Value = -19694.5005   # N136: -31249.5 -> -19694.5005

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 4827.3696   # H107: 7212.7095 -> 4827.3696
$ws.Cells.Item(107, 9).Value = 4189.0527   # I107: 6987.7827 -> 4189.0527
$ws.Cells.Item(107, 11).Value = 4189.0527   # K107: 6987.7827 -> 4189.0527
$ws.Cells.Item(107, 13).Value = -2269.0527   # M107: -5067.7827 -> -2269.0527
$ws.Cells.Item(134, 8).Value = 3022.4255   # H134: 2241.8823 -> 3022.4255
$ws.Cells.Item(134, 9).Value = 2806.4634   # I134: 2023.1587 -> 2806.4634
$ws.Cells.Item(134, 10).Value = 4498.1665   # J134: 4997.8 -> 4498.1665
$ws.Cells.Item(134, 11).Value = 8419.3902   # K134: 6069.4761 -> 8419.3902
$ws.Cells.Item(134, 12).Value = 13494.4995   # L134: 14993.4 -> 13494.4995
$ws.Cells.Item(134, 13).Value = -5884.3902   # M134: -3534.4761 -> -5884.3902
$ws.Cells.Item(134, 14).Value = -18564.4995   # N134: -20063.4 -> -18564.4995

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2856.25   # H31: 3504.4546 -> 2856.25
$ws.Cells.Item(31, 9).Value = 2257.24   # I31: 3089.6667 -> 2257.24
$ws.Cells.Item(31, 10).Value = 3339.3225   # J31: 3791.6155 -> 3339.3225
$ws.Cells.Item(31, 11).Value = 2257.24   # K31: 3089.6667 -> 2257.24
$ws.Cells.Item(31, 12).Value = 3339.3225   # L31: 3791.6155 -> 3339.3225
$ws.Cells.Item(31, 13).Value = -1962.24   # M31: -2794.6667 -> -1962.24
$ws.Cells.Item(31, 14).Value = -3929.3225   # N31: -4381.6155 -> -3929.3225
$ws.Cells.Item(34, 8).Value = 2856.25   # H34: 3504.4546 -> 2856.25
$ws.Cells.Item(34, 9).Value = 2257.24   # I34: 3089.6667 -> 2257.24
$ws.Cells.Item(34, 10).Value = 3339.3225   # J34: 3791.6155 -> 3339.3225
$ws.Cells.Item(34, 11).Value = 2257.24   # K34: 3089.6667 -> 2257.24
$ws.Cells.Item(34, 12).Value = 3339.3225   # L34: 3791.6155 -> 3339.3225
$ws.Cells.Item(34, 13).Value = -2055.24   # M34: -2887.6667 -> -2055.24
$ws.Cells.Item(34, 14).Value = -3743.3225   # N34: -4195.6155 -> -3743.3225
$ws.Cells.Item(58, 8).Value = 5461.6855   # H58: 6025.577 -> 5461.6855
$ws.Cells.Item(58, 9).Value = 5831.8945   # I58: 5729.75 -> 5831.8945
$ws.Cells.Item(58, 10).Value = 5022.0625   # J58: 6498.9 -> 5022.0625
$ws.Cells.Item(58, 11).Value = 5831.8945   # K58: 5729.75 -> 5831.8945
$ws.Cells.Item(58, 12).Value = 5022.0625   # L58: 6498.9 -> 5022.0625
$ws.Cells.Item(58, 13).Value = -5628.8945   # M58: -5526.75 -> -5628.8945
$ws.Cells.Item(58, 14).Value = -5428.0625   # N58: -6904.9 -> -5428.0625
$ws.Cells.Item(122, 8).Value = 1577.683   # H122: 1707.5946 -> 1577.683
$ws.Cells.Item(122, 9).Value = 1350.25   # I122: 1419.8 -> 1350.25
$ws.Cells.Item(122, 10).Value = 2386.3333   # J122: 2941 -> 2386.3333
$ws.Cells.Item(122, 11).Value = 4050.75   # K122: 4259.4 -> 4050.75
$ws.Cells.Item(122, 12).Value = 7158.999899999999   # L122: 8823 -> 7158.999899999999
$ws.Cells.Item(122, 13).Value = -1600.75   # M122: -1809.4 -> -1600.75
$ws.Cells.Item(122, 14).Value = -12058.9999   # N122: -13723 -> -12058.9999
$ws.Cells.Item(134, 8).Value = 5387.3096   # H134: 3327.2166 -> 5387.3096
$ws.Cells.Item(134, 9).Value = 4124.3438   # I134: 2538.9795 -> 4124.3438
$ws.Cells.Item(134, 10).Value = 9428.799999999999   # J134: 6838.4546 -> 9428.799999999999
$ws.Cells.Item(134, 11).Value = 12373.0314   # K134: 7616.9385 -> 12373.0314
$ws.Cells.Item(134, 12).Value = 28286.4   # L134: 20515.3638 -> 28286.4
$ws.Cells.Item(134, 13).Value = -9838.0314   # M134: -5081.9385 -> -9838.0314
$ws.Cells.Item(134, 14).Value = -33356.39999999999   # N134: -25585.3638 -> -33356.39999999999
$ws.Cells.Item(136, 8).Value = 5461.6855   # H136: 6025.577 -> 5461.6855
$ws.Cells.Item(136, 9).Value = 5831.8945   # I136: 5729.75 -> 5831.8945
$ws.Cells.Item(136, 10).Value = 5022.0625   # J136: 6498.9 -> 5022.0625
$ws.Cells.Item(136, 11).Value = 17495.6835   # K136: 17189.25 -> 17495.6835
$ws.Cells.Item(136, 12).Value = 15066.1875   # L136: 19496.7 -> 15066.1875
$ws.Cells.Item(136, 13).Value = -14945.6835   # M136: -14639.25 -> -14945.6835
$ws.Cells.Item(136, 14).Value = -20166.1875   # N136: -24596.7 -> -20166.1875
$ws.Cells.Item(138, 8).Value = 83406.39999999999   # H138: 83732.3 -> 83406.39999999999
$ws.Cells.Item(138, 10).Value = 83406.39999999999   # J138: 83732.3 -> 83406.39999999999
$ws.Cells.Item(138, 12).Value = 83406.39999999999   # L138: 83732.3 -> 83406.39999999999
$ws.Cells.Item(138, 14).Value = -93686.39999999999   # N138: -94012.3 -> -93686.39999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 652.1429000000001   # H5: 800.8889 -> 652.1429000000001
$ws.Cells.Item(5, 9).Value = 703.7778   # I5: 800.8889 -> 703.7778
$ws.Cells.Item(5, 10).Value = 559.2   # J5: 0 -> 559.2
$ws.Cells.Item(5, 11).Value = 2111.3334   # K5: 2402.6667 -> 2111.3334
$ws.Cells.Item(5, 12).Value = 1677.6   # L5: 0 -> 1677.6
$ws.Cells.Item(5, 13).Value = -1999.3334   # M5: -2290.6667 -> -1999.3334
$ws.Cells.Item(5, 14).Value = -1901.6   # N5: None -> -1901.6
$ws.Cells.Item(11, 8).Value = 63512.875   # H11: 63664 -> 63512.875
$ws.Cells.Item(11, 9).Value = 91469.55   # I11: 100358.5 -> 91469.55
$ws.Cells.Item(11, 10).Value = 2008.2   # J11: 2506.5 -> 2008.2
$ws.Cells.Item(11, 11).Value = 274408.65   # K11: 301075.5 -> 274408.65
$ws.Cells.Item(11, 12).Value = 6024.6   # L11: 7519.5 -> 6024.6
$ws.Cells.Item(11, 13).Value = -274268.65   # M11: -300935.5 -> -274268.65
$ws.Cells.Item(11, 14).Value = -6304.6   # N11: -7799.5 -> -6304.6
$ws.Cells.Item(68, 8).Value = 2566.0303   # H68: 2582.9395 -> 2566.0303
$ws.Cells.Item(68, 9).Value = 2300.1428   # I68: 2262.625 -> 2300.1428
$ws.Cells.Item(68, 10).Value = 2637.6155   # J68: 2685.44 -> 2637.6155
$ws.Cells.Item(68, 11).Value = 6900.428400000001   # K68: 6787.875 -> 6900.428400000001
$ws.Cells.Item(68, 12).Value = 7912.8465   # L68: 8056.32 -> 7912.8465
$ws.Cells.Item(68, 13).Value = -6089.428400000001   # M68: -5976.875 -> -6089.428400000001
$ws.Cells.Item(68, 14).Value = -9534.8465   # N68: -9678.32 -> -9534.8465
$ws.Cells.Item(71, 8).Value = 2566.0303   # H71: 2582.9395 -> 2566.0303
$ws.Cells.Item(71, 9).Value = 2300.1428   # I71: 2262.625 -> 2300.1428
$ws.Cells.Item(71, 10).Value = 2637.6155   # J71: 2685.44 -> 2637.6155
$ws.Cells.Item(71, 11).Value = 20701.2852   # K71: 20363.625 -> 20701.2852
$ws.Cells.Item(71, 12).Value = 23738.5395   # L71: 24168.96 -> 23738.5395
$ws.Cells.Item(71, 13).Value = -16645.2852   # M71: -16307.625 -> -16645.2852
$ws.Cells.Item(71, 14).Value = -31850.5395   # N71: -32280.96 -> -31850.5395
$ws.Cells.Item(113, 8).Value = 548.9231   # H113: 3234.25 -> 548.9231
$ws.Cells.Item(113, 9).Value = 198.6   # I113: 1899.5 -> 198.6
$ws.Cells.Item(113, 10).Value = 767.875   # J113: 4569 -> 767.875
$ws.Cells.Item(113, 11).Value = 595.8   # K113: 5698.5 -> 595.8
$ws.Cells.Item(113, 12).Value = 2303.625   # L113: 13707 -> 2303.625
$ws.Cells.Item(113, 13).Value = 1574.2   # M113: -3528.5 -> 1574.2
$ws.Cells.Item(113, 14).Value = -6643.625   # N113: -18047 -> -6643.625
$ws.Cells.Item(131, 8).Value = 6811932.5   # H131: 7018333.5 -> 6811932.5
$ws.Cells.Item(131, 9).Value = 15041255   # I131: 17189906 -> 15041255
$ws.Cells.Item(131, 11).Value = 45123765   # K131: 51569718 -> 45123765
$ws.Cells.Item(131, 13).Value = -45118725   # M131: -51564678 -> -45118725
$ws.Cells.Item(133, 8).Value = 1999   # H133: 2981.6667 -> 1999
$ws.Cells.Item(133, 9).Value = 1999   # I133: 2981.6667 -> 1999
$ws.Cells.Item(133, 11).Value = 5997   # K133: 8945.000100000001 -> 5997
$ws.Cells.Item(133, 13).Value = -937   # M133: -3885.000100000001 -> -937
$ws.Cells.Item(134, 8).Value = 2901.6155   # H134: 2020.1111 -> 2901.6155
$ws.Cells.Item(134, 9).Value = 2047.1   # I134: 2020.1111 -> 2047.1
$ws.Cells.Item(134, 10).Value = 5750   # J134: 0 -> 5750
$ws.Cells.Item(134, 11).Value = 6141.299999999999   # K134: 6060.3333 -> 6141.299999999999
$ws.Cells.Item(134, 12).Value = 17250   # L134: 0 -> 17250
$ws.Cells.Item(134, 13).Value = -1071.299999999999   # M134: -990.3333000000002 -> -1071.299999999999
$ws.Cells.Item(134, 14).Value = -27390   # N134: None -> -27390
$ws.Cells.Item(135, 8).Value = 652.1429000000001   # H135: 800.8889 -> 652.1429000000001
$ws.Cells.Item(135, 9).Value = 703.7778   # I135: 800.8889 -> 703.7778
$ws.Cells.Item(135, 10).Value = 559.2   # J135: 0 -> 559.2
$ws.Cells.Item(135, 11).Value = 6334.000199999999   # K135: 7208.0001 -> 6334.000199999999
$ws.Cells.Item(135, 12).Value = 5032.8   # L135: 0 -> 5032.8
$ws.Cells.Item(135, 13).Value = -3799.000199999999   # M135: -4673.0001 -> -3799.000199999999
$ws.Cells.Item(135, 14).Value = -10102.8   # N135: None -> -10102.8
$ws.Cells.Item(140, 8).Value = 672.6896400000001   # H140: 891.9 -> 672.6896400000001
$ws.Cells.Item(140, 9).Value = 681.8148   # I140: 896.7368 -> 681.8148
$ws.Cells.Item(140, 10).Value = 549.5   # J140: 800 -> 549.5
$ws.Cells.Item(140, 11).Value = 2045.4444   # K140: 2690.2104 -> 2045.4444
$ws.Cells.Item(140, 12).Value = 1648.5   # L140: 2400 -> 1648.5
$ws.Cells.Item(140, 13).Value = 3134.5556   # M140: 2489.7896 -> 3134.5556
$ws.Cells.Item(140, 14).Value = -12008.5   # N140: -12760 -> -12008.5
$ws.Cells.Item(141, 8).Value = 21927.2   # H141: 17838 -> 21927.2
$ws.Cells.Item(141, 9).Value = 21927.2   # I141: 17838 -> 21927.2
$ws.Cells.Item(141, 11).Value = 65781.60000000001   # K141: 53514 -> 65781.60000000001
$ws.Cells.Item(141, 13).Value = -60601.60000000001   # M141: -48334 -> -60601.60000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(92, 8).Value = 30483.334   # H92: 34725 -> 30483.334
$ws.Cells.Item(92, 10).Value = 30483.334   # J92: 34725 -> 30483.334
$ws.Cells.Item(92, 12).Value = 30483.334   # L92: 34725 -> 30483.334
$ws.Cells.Item(92, 14).Value = -34227.334   # N92: -38469 -> -34227.334
$ws.Cells.Item(122, 8).Value = 2323.4167   # H122: 1788.2632 -> 2323.4167
$ws.Cells.Item(122, 9).Value = 2515.0908   # I122: 1554.84 -> 2515.0908
$ws.Cells.Item(122, 10).Value = 2161.2307   # J122: 2237.1538 -> 2161.2307
$ws.Cells.Item(122, 11).Value = 7545.2724   # K122: 4664.52 -> 7545.2724
$ws.Cells.Item(122, 12).Value = 6483.6921   # L122: 6711.4614 -> 6483.6921
$ws.Cells.Item(122, 13).Value = -5095.2724   # M122: -2214.52 -> -5095.2724
$ws.Cells.Item(122, 14).Value = -11383.6921   # N122: -11611.4614 -> -11383.6921
$ws.Cells.Item(140, 8).Value = 69355.625   # H140: 109573.75 -> 69355.625
$ws.Cells.Item(140, 10).Value = 69355.625   # J140: 109573.75 -> 69355.625
$ws.Cells.Item(140, 12).Value = 69355.625   # L140: 109573.75 -> 69355.625
$ws.Cells.Item(140, 14).Value = -79715.625   # N140: -119933.75 -> -79715.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 2438.1667   # H16: 3414.0557 -> 2438.1667
$ws.Cells.Item(16, 9).Value = 2417.9375   # I16: 3113.5386 -> 2417.9375
$ws.Cells.Item(16, 10).Value = 2600   # J16: 4195.4 -> 2600
$ws.Cells.Item(16, 11).Value = 2417.9375   # K16: 3113.5386 -> 2417.9375
$ws.Cells.Item(16, 12).Value = 2600   # L16: 4195.4 -> 2600
$ws.Cells.Item(16, 13).Value = -2247.9375   # M16: -2943.5386 -> -2247.9375
$ws.Cells.Item(16, 14).Value = -2940   # N16: -4535.4 -> -2940
$ws.Cells.Item(55, 8).Value = 583.3   # H55: 1016.9091 -> 583.3
$ws.Cells.Item(55, 9).Value = 176.33333   # I55: 380.36365 -> 176.33333
$ws.Cells.Item(55, 10).Value = 916.2727   # J55: 1653.4546 -> 916.2727
$ws.Cells.Item(55, 11).Value = 176.33333   # K55: 380.36365 -> 176.33333
$ws.Cells.Item(55, 12).Value = 916.2727   # L55: 1653.4546 -> 916.2727
$ws.Cells.Item(55, 13).Value = -3.333329999999989   # M55: -207.36365 -> -3.333329999999989
$ws.Cells.Item(55, 14).Value = -1262.2727   # N55: -1999.4546 -> -1262.2727
$ws.Cells.Item(108, 8).Value = 9313   # H108: 0 -> 9313
$ws.Cells.Item(108, 10).Value = 9313   # J108: 0 -> 9313
$ws.Cells.Item(108, 12).Value = 9313   # L108: 0 -> 9313
$ws.Cells.Item(108, 14).Value = -16993   # N108: None -> -16993
$ws.Cells.Item(122, 8).Value = 4378.4375   # H122: 5173.1665 -> 4378.4375
$ws.Cells.Item(122, 9).Value = 3037.25   # I122: 3584.125 -> 3037.25
$ws.Cells.Item(122, 10).Value = 8402   # J122: 8351.25 -> 8402
$ws.Cells.Item(122, 11).Value = 9111.75   # K122: 10752.375 -> 9111.75
$ws.Cells.Item(122, 12).Value = 25206   # L122: 25053.75 -> 25206
$ws.Cells.Item(122, 13).Value = -6661.75   # M122: -8302.375 -> -6661.75
$ws.Cells.Item(122, 14).Value = -30106   # N122: -29953.75 -> -30106
$ws.Cells.Item(132, 8).Value = 58299.367   # H132: 48703.78 -> 58299.367
$ws.Cells.Item(132, 9).Value = 72112.8   # I132: 63987.824 -> 72112.8
$ws.Cells.Item(132, 10).Value = 6499   # J132: 5399 -> 6499
$ws.Cells.Item(132, 11).Value = 216338.4   # K132: 191963.472 -> 216338.4
$ws.Cells.Item(132, 12).Value = 19497   # L132: 16197 -> 19497
$ws.Cells.Item(132, 13).Value = -213808.4   # M132: -189433.472 -> -213808.4
$ws.Cells.Item(132, 14).Value = -24557   # N132: -21257 -> -24557

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 1819.8   # H132: 1572.9375 -> 1819.8
$ws.Cells.Item(132, 9).Value = 1524.75   # I132: 1315.2 -> 1524.75
$ws.Cells.Item(132, 10).Value = 3000   # J132: 2493.4285 -> 3000
$ws.Cells.Item(132, 11).Value = 4574.25   # K132: 3945.6 -> 4574.25
$ws.Cells.Item(132, 12).Value = 9000   # L132: 7480.2855 -> 9000
$ws.Cells.Item(132, 13).Value = -2044.25   # M132: -1415.6 -> -2044.25
$ws.Cells.Item(132, 14).Value = -14060   # N132: -12540.2855 -> -14060
$ws.Cells.Item(136, 8).Value = 35248.8   # H136: 32829.37 -> 35248.8
$ws.Cells.Item(136, 9).Value = 40874.445   # I136: 35471.81 -> 40874.445
$ws.Cells.Item(136, 10).Value = 20782.857   # J136: 23580.834 -> 20782.857
$ws.Cells.Item(136, 11).Value = 122623.335   # K136: 106415.43 -> 122623.335
$ws.Cells.Item(136, 12).Value = 62348.571   # L136: 70742.50199999999 -> 62348.571
$ws.Cells.Item(136, 13).Value = -120073.335   # M136: -103865.43 -> -120073.335
$ws.Cells.Item(136, 14).Value = -67448.571   # N136: -75842.50199999999 -> -67448.571
$ws.Cells.Item(139, 8).Value = 0   # H139: 69804.86 -> 0
$ws.Cells.Item(139, 9).Value = 0   # I139: 69000 -> 0
$ws.Cells.Item(139, 10).Value = 0   # J139: 69939 -> 0
$ws.Cells.Item(139, 11).Value = 0   # K139: 69000 -> 0
$ws.Cells.Item(139, 12).Value = 0   # L139: -80219 -> 0
$ws.Cells.Item(139, 13).ClearContents()   # M139: removed
$ws.Cells.Item(139, 14).ClearContents()   # N139: removed
